$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: A25 becomes a genuine numeric value (was stored as text before)
$ws.Range("A25").Value = 71277628

# New row 26 - payment 71277628 (Cash) 2025-08-18T16:54:45
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "71277628"
$ws.Range("A26").ClearFormats()

$ws.Range("B26").Value = ""
$ws.Range("C26").Value = "Cash"
$ws.Range("D26").Value = "2025-08-18T16:54:45"
$ws.Range("E26").Value = 76
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 76
